$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = -7.343399999999995
$ws.Range("D18").Value = -8.907199999999996
$ws.Range("D20").Value = -7.6062
$ws.Range("D27").Value = -9.014199999999995
$ws.Range("D35").Value = -8.070400000000001
$ws.Range("D69").Value = -7.072699999999998
$ws.Range("D76").Value = -7.444200000000002
$ws.Range("D78").Value = -7.570200000000005
$ws.Range("D82").Value = -8.315299999999993
$ws.Range("D83").Value = -9.0831
$ws.Range("D93").Value = -6.503199999999997
